$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text format first so numeric-looking values
# ("1.001", "42.15", etc.) are stored as text, matching the source data
# which uses inline strings, not numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '22.015.70'
$ws.Range("E2").Value = '  -1.51%  '
$ws.Range("D3").Value = '1.551.93'
$ws.Range("E3").Value = '  -0.80%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("E5").Value = '  +0.00%  '
$ws.Range("D6").Value = '287.51'
$ws.Range("E6").Value = '  +0.45%  '
$ws.Range("D7").Value = '0.3928'
$ws.Range("E7").Value = '  +4.82%  '
$ws.Range("D8").Value = '0.3195'
$ws.Range("E8").Value = '  -2.20%  '
$ws.Range("D9").Value = '42.15'
$ws.Range("E9").Value = '  -7.12%  '
$ws.Range("D10").Value = '0.07264'
$ws.Range("E10").Value = '  -1.73%  '
$ws.Range("D11").Value = '1.090'
$ws.Range("E11").Value = '  -4.32%  '
$ws.Range("D12").Value = '1.001'
$ws.Range("E12").Value = '  -0.07%  '
$ws.Range("D13").Value = '18.88'
$ws.Range("E13").Value = '  -7.22%  '
$ws.Range("D14").Value = '5.602'
$ws.Range("E14").Value = '  -3.80%  '
$ws.Range("D15").Value = '6.651'
$ws.Range("E15").Value = '  -2.16%  '
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").Value = '0.00001119'
$ws.Range("E16").Value = '  +2.28%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '1.549.49'
$ws.Range("E17").Value = '  -0.87%  '
$ws.Range("D18").Value = '0.06582'
$ws.Range("E18").Value = '  -1.96%  '
$ws.Range("D19").Value = '83.81'
$ws.Range("E19").Value = '  -2.29%  '
$ws.Range("E20").Value = '  +0.02%  '
$ws.Range("D21").Value = '6.284'
$ws.Range("E21").Value = '  -0.78%  '
$ws.Range("D22").Value = '15.70'
$ws.Range("E22").Value = '  -3.28%  '
$ws.Range("D23").Value = '11.19'
$ws.Range("E23").Value = '  -4.00%  '
$ws.Range("D24").Value = '22.015.80'
$ws.Range("E24").Value = '  -1.54%  '
$ws.Range("D25").Value = '2.348'
$ws.Range("D26").Value = '2.426'
$ws.Range("E26").Value = '  -3.74%  '
$ws.Range("D27").Value = '147.08'
$ws.Range("E27").Value = '  -1.77%  '
$ws.Range("D28").Value = '18.59'
$ws.Range("E28").Value = '  -4.05%  '
$ws.Range("D29").Value = '4.834'
$ws.Range("D30").Value = '1.725.24'
$ws.Range("E30").Value = '  -0.84%  '
$ws.Range("D31").Value = '119.05'
$ws.Range("E31").Value = '  -3.16%  '
$ws.Range("D32").Value = '1.062'
$ws.Range("E32").Value = '  +1.03%  '
$ws.Range("D33").Value = '5.660'
$ws.Range("E33").Value = '  -4.04%  '
$ws.Range("D34").Value = '0.08313'
$ws.Range("E34").Value = '  +1.39%  '
$ws.Range("D35").Value = '9.169'
$ws.Range("E35").Value = '  -3.39%  '
$ws.Range("D36").Value = '1.592'
$ws.Range("E36").Value = '  -16.80%  '
$ws.Range("D37").Value = '0.06155'
$ws.Range("E37").Value = '  -2.08%  '
$ws.Range("D38").Value = '0.02256'
$ws.Range("E38").Value = '  -5.14%  '
$ws.Range("D39").Value = '5.091'
$ws.Range("E39").Value = '  -2.93%  '
$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").Value = '1.212'
$ws.Range("E40").Value = '  -5.60%  '
$ws.Range("B41").Value = 'Algorand'
$ws.Range("C41").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D41").Value = '0.2061'
$ws.Range("E41").Value = '  -5.29%  '
$ws.Range("E42").Value = '  +0.00%  '
$ws.Range("D43").Value = '10.54'
$ws.Range("E43").Value = '  -4.06%  '
$ws.Range("D44").Value = '0.5790'
$ws.Range("E44").Value = '  -4.68%  '
$ws.Range("D45").Value = '13.16'
$ws.Range("E45").Value = '  -4.21%  '
$ws.Range("D46").Value = '3.706'
$ws.Range("E46").Value = '  -0.86%  '
$ws.Range("D47").Value = '0.5534'
$ws.Range("E47").Value = '  -6.02%  '
$ws.Range("D48").Value = '1.892'
$ws.Range("E48").Value = '  -4.96%  '
$ws.Range("D49").Value = '117.24'
$ws.Range("E49").Value = '  -5.16%  '
$ws.Range("D50").Value = '1.136'
$ws.Range("E50").Value = '  -3.44%  '
$ws.Range("D51").Value = '0.06827'
$ws.Range("E51").Value = '  -4.27%  '

# Restore the default (Normal) style on column D so no stray number
# format / style index is left behind on the cells.
$ws.Range("D2:D51").Style = "Normal"